# Season-record columns (Wins / Losses / Ties) were missing from the
# original scrape. Add them as the three new trailing columns AD:AF,
# mirroring the header style used by the existing header row and
# filling every data row (2-46) with the team's season record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy the formatting of an existing header cell (bold,
# centered, thin border - style index 1) onto AD1:AF1, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row shares the same team season record: 89 wins, 73 losses,
# 0 ties.
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 73
    $ws.Cells.Item($r, 32).Value = 0
}
